{"js": "// The document contains a single 20x5 table of simple arithmetic\n// expressions (e.g. \"38+7=\"). The commit replaces the text of 100\n// cells (in row-major order) with a new set of expressions, leaving\n// every other part of the document (fonts, sizes, paragraph\n// alignment, table/row/cell structure) untouched.\n//\n// Because a couple of the original expressions are not unique\n// (e.g. \"11+58=\" appears twice but maps to two different\n// replacements), the mapping below is applied strictly by position\n// (row-major cell order), not by text search/replace.\n\nconst beforeRows = [\n  [\"38+7=\", \"18-10=\", \"89-7=\", \"48+14=\", \"62+33=\"],\n  [\"23+29=\", \"55+35=\", \"23+8=\", \"71-47=\", \"2+27=\"],\n  [\"10+19=\", \"25+36=\", \"82-2=\", \"90-31=\", \"96-16=\"],\n  [\"40+35=\", \"69-13=\", \"89-61=\", \"58-3=\", \"83-45=\"],\n  [\"1+74=\", \"23+74=\", \"35+17=\", \"18-0=\", \"27+3=\"],\n  [\"20+19=\", \"38+26=\", \"90-9=\", \"92-52=\", \"7+3=\"],\n  [\"11+71=\", \"69-15=\", \"3+38=\", \"68-7=\", \"74-43=\"],\n  [\"43+48=\", \"10+55=\", \"87-79=\", \"57+11=\", \"34+8=\"],\n  [\"67-45=\", \"15+4=\", \"49+50=\", \"8+40=\", \"13+60=\"],\n  [\"94-2=\", \"71-14=\", \"46+15=\", \"76-73=\", \"77-74=\"],\n  [\"88-28=\", \"97-8=\", \"47-5=\", \"79-33=\", \"10+24=\"],\n  [\"67-32=\", \"52+36=\", \"73-43=\", \"6+73=\", \"61-6=\"],\n  [\"45+50=\", \"41-8=\", \"1+70=\", \"52-41=\", \"43+8=\"],\n  [\"77-48=\", \"22-8=\", \"32+33=\", \"34-19=\", \"21-8=\"],\n  [\"20+16=\", \"35+25=\", \"31-5=\", \"64-16=\", \"5+27=\"],\n  [\"76-58=\", \"11+58=\", \"31-16=\", \"97-89=\", \"11+44=\"],\n  [\"90-14=\", \"63-41=\", \"83-58=\", \"11+58=\", \"43-7=\"],\n  [\"5+15=\", \"54+32=\", \"76-53=\", \"79-3=\", \"10-5=\"],\n  [\"18+76=\", \"73+18=\", \"25+14=\", \"94-18=\", \"78-73=\"],\n  [\"99-27=\", \"28+54=\", \"79-15=\", \"25-15=\", \"28+31=\"],\n];\n\nconst afterRows = [\n  [\"79-28=\", \"34-9=\", \"28-25=\", \"87-37=\", \"39+27=\"],\n  [\"99-50=\", \"80+18=\", \"79-18=\", \"73-0=\", \"17+23=\"],\n  [\"12+81=\", \"87-56=\", \"41+7=\", \"16+43=\", \"97-90=\"],\n  [\"20+25=\", \"64-33=\", \"52-39=\", \"37+48=\", \"85-2=\"],\n  [\"49+27=\", \"48-38=\", \"30+19=\", \"95-35=\", \"64-6=\"],\n  [\"35+34=\", \"33-32=\", \"99-98=\", \"8+89=\", \"93-44=\"],\n  [\"4+59=\", \"76-67=\", \"70+4=\", \"41+48=\", \"54+30=\"],\n  [\"80-37=\", \"29-28=\", \"76-4=\", \"78-65=\", \"52-6=\"],\n  [\"4+71=\", \"42-26=\", \"31-24=\", \"94-1=\", \"97-93=\"],\n  [\"17-8=\", \"28-4=\", \"28-22=\", \"70-57=\", \"94-35=\"],\n  [\"76-1=\", \"61+20=\", \"19+12=\", \"94-73=\", \"69-21=\"],\n  [\"69-11=\", \"20+24=\", \"87-28=\", \"76-57=\", \"55+1=\"],\n  [\"46-45=\", \"30+22=\", \"44+38=\", \"71+24=\", \"12+66=\"],\n  [\"18+51=\", \"63-18=\", \"52+19=\", \"95-64=\", \"25+23=\"],\n  [\"2+70=\", \"70-18=\", \"96-67=\", \"18+20=\", \"56+17=\"],\n  [\"32+27=\", \"57-20=\", \"55-10=\", \"6+32=\", \"63+28=\"],\n  [\"25+24=\", \"68-48=\", \"59+23=\", \"69-31=\", \"52+21=\"],\n  [\"38+16=\", \"69-66=\", \"70-27=\", \"25+56=\", \"40+3=\"],\n  [\"16+68=\", \"78-37=\", \"38+1=\", \"95-4=\", \"53-6=\"],\n  [\"42+12=\", \"74-66=\", \"66+4=\", \"18+40=\", \"87-85=\"],\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"values,rowCount\");\nawait context.sync();\n\nconst current = table.values;\n\n// Build the replacement grid: start from whatever is currently in the\n// table, and only overwrite a cell when it still holds the expected\n// \"before\" value at that position (row-major). This keeps the edit\n// anchored to the documented diff even if unrelated cells were\n// touched by something else, while still coping with duplicate\n// \"before\" text values because matching is positional.\nconst next = current.map((row, r) =>\n  row.map((text, c) => {\n    const expectedBefore = beforeRows[r] && beforeRows[r][c];\n    const expectedAfter = afterRows[r] && afterRows[r][c];\n    if (expectedAfter === undefined) {\n      return text;\n    }\n    if (text === expectedBefore || expectedBefore === undefined) {\n      return expectedAfter;\n    }\n    return text;\n  })\n);\n\ntable.values = next;\nawait context.sync();\n", "ps1": "# The document contains a single 20x5 table of simple arithmetic\n# expressions (e.g. \"38+7=\"). The commit replaces the text of 100\n# cells (in row-major order) with a new set of expressions, leaving\n# every other part of the document (fonts, sizes, paragraph\n# alignment, table/row/cell structure) untouched.\n#\n# Because a couple of the original expressions are not unique\n# (e.g. \"11+58=\" appears twice but maps to two different\n# replacements), the mapping below is applied strictly by position\n# (row index, column index), not by text search/replace.\n\n$beforeRows = @(\n    @(\"38+7=\", \"18-10=\", \"89-7=\", \"48+14=\", \"62+33=\"),\n    @(\"23+29=\", \"55+35=\", \"23+8=\", \"71-47=\", \"2+27=\"),\n    @(\"10+19=\", \"25+36=\", \"82-2=\", \"90-31=\", \"96-16=\"),\n    @(\"40+35=\", \"69-13=\", \"89-61=\", \"58-3=\", \"83-45=\"),\n    @(\"1+74=\", \"23+74=\", \"35+17=\", \"18-0=\", \"27+3=\"),\n    @(\"20+19=\", \"38+26=\", \"90-9=\", \"92-52=\", \"7+3=\"),\n    @(\"11+71=\", \"69-15=\", \"3+38=\", \"68-7=\", \"74-43=\"),\n    @(\"43+48=\", \"10+55=\", \"87-79=\", \"57+11=\", \"34+8=\"),\n    @(\"67-45=\", \"15+4=\", \"49+50=\", \"8+40=\", \"13+60=\"),\n    @(\"94-2=\", \"71-14=\", \"46+15=\", \"76-73=\", \"77-74=\"),\n    @(\"88-28=\", \"97-8=\", \"47-5=\", \"79-33=\", \"10+24=\"),\n    @(\"67-32=\", \"52+36=\", \"73-43=\", \"6+73=\", \"61-6=\"),\n    @(\"45+50=\", \"41-8=\", \"1+70=\", \"52-41=\", \"43+8=\"),\n    @(\"77-48=\", \"22-8=\", \"32+33=\", \"34-19=\", \"21-8=\"),\n    @(\"20+16=\", \"35+25=\", \"31-5=\", \"64-16=\", \"5+27=\"),\n    @(\"76-58=\", \"11+58=\", \"31-16=\", \"97-89=\", \"11+44=\"),\n    @(\"90-14=\", \"63-41=\", \"83-58=\", \"11+58=\", \"43-7=\"),\n    @(\"5+15=\", \"54+32=\", \"76-53=\", \"79-3=\", \"10-5=\"),\n    @(\"18+76=\", \"73+18=\", \"25+14=\", \"94-18=\", \"78-73=\"),\n    @(\"99-27=\", \"28+54=\", \"79-15=\", \"25-15=\", \"28+31=\")\n)\n\n$afterRows = @(\n    @(\"79-28=\", \"34-9=\", \"28-25=\", \"87-37=\", \"39+27=\"),\n    @(\"99-50=\", \"80+18=\", \"79-18=\", \"73-0=\", \"17+23=\"),\n    @(\"12+81=\", \"87-56=\", \"41+7=\", \"16+43=\", \"97-90=\"),\n    @(\"20+25=\", \"64-33=\", \"52-39=\", \"37+48=\", \"85-2=\"),\n    @(\"49+27=\", \"48-38=\", \"30+19=\", \"95-35=\", \"64-6=\"),\n    @(\"35+34=\", \"33-32=\", \"99-98=\", \"8+89=\", \"93-44=\"),\n    @(\"4+59=\", \"76-67=\", \"70+4=\", \"41+48=\", \"54+30=\"),\n    @(\"80-37=\", \"29-28=\", \"76-4=\", \"78-65=\", \"52-6=\"),\n    @(\"4+71=\", \"42-26=\", \"31-24=\", \"94-1=\", \"97-93=\"),\n    @(\"17-8=\", \"28-4=\", \"28-22=\", \"70-57=\", \"94-35=\"),\n    @(\"76-1=\", \"61+20=\", \"19+12=\", \"94-73=\", \"69-21=\"),\n    @(\"69-11=\", \"20+24=\", \"87-28=\", \"76-57=\", \"55+1=\"),\n    @(\"46-45=\", \"30+22=\", \"44+38=\", \"71+24=\", \"12+66=\"),\n    @(\"18+51=\", \"63-18=\", \"52+19=\", \"95-64=\", \"25+23=\"),\n    @(\"2+70=\", \"70-18=\", \"96-67=\", \"18+20=\", \"56+17=\"),\n    @(\"32+27=\", \"57-20=\", \"55-10=\", \"6+32=\", \"63+28=\"),\n    @(\"25+24=\", \"68-48=\", \"59+23=\", \"69-31=\", \"52+21=\"),\n    @(\"38+16=\", \"69-66=\", \"70-27=\", \"25+56=\", \"40+3=\"),\n    @(\"16+68=\", \"78-37=\", \"38+1=\", \"95-4=\", \"53-6=\"),\n    @(\"42+12=\", \"74-66=\", \"66+4=\", \"18+40=\", \"87-85=\")\n)\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\nfor ($r = 0; $r -lt $beforeRows.Length; $r++) {\n    $beforeRow = $beforeRows[$r]\n    $afterRow = $afterRows[$r]\n    for ($c = 0; $c -lt $beforeRow.Length; $c++) {\n        $cell = $t.Cell($r + 1, $c + 1)\n        $current = $cell.Range.Text.TrimEnd([char]13, [char]7)\n        $expectedBefore = $beforeRow[$c]\n        $expectedAfter = $afterRow[$c]\n        if ($current -eq $expectedBefore -or $expectedBefore -eq $null) {\n            $cell.Range.Text = $expectedAfter\n        }\n    }\n}\n"}
